# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "conversion" text block in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n" +
           "✅ Dólar paralelo: 68`n" +
           "`n" +
           "Binance`n" +
           "✅ 1000 Bs = 7.04 = 28211.27 pesos`n" +
           "✅ 28211.27 pesos = 7.02 = 978.27 Bs`n" +
           "`n" +
           "Promedio competencia`n" +
           "✅ Tasa pesos: 20`n" +
           "✅ Tasa Bs: 20`n" +
           "✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- tasas: update N10, O10, N12, O12 values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 142
$wsTasas.Range("O10").Value = 4006
$wsTasas.Range("N12").Value = 4020
$wsTasas.Range("O12").Value = 139.4
